# Atualiza dados em 28-12-2017
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - FHEMIG
$ws.Range("E2").Value = 1259966447
$ws.Range("F2").Value = 1251705752.77
$ws.Range("G2").Value = 1250405684.67
$ws.Range("O2").Value = 1253063721.46
$ws.Range("P2").Value = 1253063721.46

# Row 3 - SEAP
$ws.Range("E3").Value = 139884594.38
$ws.Range("F3").Value = 135676732.65
$ws.Range("G3").Value = 134335565.37
$ws.Range("O3").Value = 132741835

# Row 4 - SESP
$ws.Range("F4").Value = 12879984.11
$ws.Range("G4").Value = 11227561.51

# Row 5 - SETOP
$ws.Range("E5").Value = 7385707.129999999
$ws.Range("F5").Value = 4714403.09
$ws.Range("K5").Value = 0.000000001862645149230957

# Row 6 - UNIMONTES
$ws.Range("E6").Value = 89436899
$ws.Range("F6").Value = 89274115.61
$ws.Range("G6").Value = 89206993.31
$ws.Range("O6").Value = 0

# Row 7 - FUNED
$ws.Range("E7").Value = 398183440.28
$ws.Range("F7").Value = 394517627.99
$ws.Range("G7").Value = 349080701.49

# Row 8 - HEMOMINAS
$ws.Range("E8").Value = 246642653.08
$ws.Range("F8").Value = 243721433.44
$ws.Range("G8").Value = 242640807.93

# Row 9 - ESP-MG
$ws.Range("E9").Value = 12246442.6
$ws.Range("F9").Value = 12067121.22
$ws.Range("G9").Value = 11863643.63
$ws.Range("K9").Value = 11118132.76
$ws.Range("Q9").Value = 11118132.76

# Row 10 - FAPEMIG
$ws.Range("G10").Value = 2852691.27
$ws.Range("K10").Value = 2917300
